# Update the Modal Shift Vector values in column B as per the revised
# model_shift_matrices parameters.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 65
$ws.Range("B3").Value = 14
$ws.Range("B4").Value = 113
$ws.Range("B5").Value = 132
$ws.Range("B9").Value = 1
$ws.Range("B11").Value = 13
$ws.Range("B12").Value = 4
$ws.Range("B13").Value = 10
